# Apply the "Query" sheet changes: insert a new column B, populate new rows
# 10/11 with new SQL/description content, resize row heights, and update the
# selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "Query" sheet

# ---------------------------------------------------------------------------
# 1. Insert a new (blank) column before the current column B. This shifts the
#    existing B/C/D columns to C/D/E, matching the target layout.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ColumnWidth = 18.3

# The insert operation stamps the newly-blank B1:B7 cells with column A's
# style (since they sit next to populated A cells). The source workbook
# keeps that column genuinely empty there, so strip those phantom cells.
$ws.Range("B1:B7").Clear()

# ---------------------------------------------------------------------------
# 2. Build the new multi-line SQL strings (use CRLF to match the workbook's
#    existing convention for embedded line breaks).
# ---------------------------------------------------------------------------
$s50 = 'Get the Investigation list by prescription and visit'

$s51Lines = @(
  'SELECT b.investigation_name, a.value, b.unit, investigation_id',
  '                    FROM patient_investigation a, investigation_master b',
  '                    WHERE a.patient_id = ''123''',
  '                    AND a.visit_id = ''10404''',
  '                    AND a.investigation_id = b.ID and a.chamber_id=b.chamber_id and a.doc_id=b.doc_id and',
  '     a.chamber_id=''sos'' AND a.doc_id=''sroy'''
)
$s51 = [string]::Join("`r`n", $s51Lines)

$s52Lines = @(
  'SELECT b.investigation_name, a.value, b.unit, investigation_id',
  '            FROM patient_investigation a, investigation_master b',
  '            WHERE a.patient_id = ''$patient_id''',
  '            AND a.visit_id = ''$visit_id''',
  '            AND a.investigation_id = b.ID and a.chamber_id=b.chamber_id and a.doc_id=b.doc_id and',
  '            a.chamber_id=''".$chamber_name."'' AND a.doc_id=''".$doc_name."'''
)
$s52 = [string]::Join("`r`n", $s52Lines)

$s53Lines = @(
  'SELECT b.type, b.ID',
  '                        FROM prescribed_cf a, clinical_impression b',
  '                        WHERE a.clinical_impression_id = b.id and a.chamber_id=b.chamber_id and a.doc_id=b.doc_id and',
  '                        AND a.prescription_id = ''$PRESCRIPTION_ID'' and a.chamber_id=''".$chamber_name."'' AND a.doc_id=''".$doc_name."'''
)
$s53 = [string]::Join("`r`n", $s53Lines)

$s54 = 'Clinical impression List by prescription'
$s55 = 'makeprescription/clinical_impression.php'

# ---------------------------------------------------------------------------
# 3. Helper to stamp a cell with the shared "Trebuchet MS, 9pt" look used
#    throughout this sheet; $wrap toggles the word-wrap style variant.
# ---------------------------------------------------------------------------
function Set-QueryCell($addr, $value, $wrap) {
    $rng = $ws.Range($addr)
    $rng.Value = $value
    $rng.Font.Name = "Trebuchet MS"
    $rng.Font.Size = 9
    if ($wrap) {
        $rng.WrapText = $true
    }
}

# ---------------------------------------------------------------------------
# 4. Row 10 - "Get the Investigation list by prescription and visit"
#    (populated in the same D -> E -> C order the original authoring used,
#    so new shared-string entries line up with the source workbook)
# ---------------------------------------------------------------------------
Set-QueryCell "D10" $s50 $false
Set-QueryCell "E10" $s51 $true
Set-QueryCell "C10" $s52 $true
$ws.Rows.Item(10).RowHeight = 105

# ---------------------------------------------------------------------------
# 5. Row 11 - "Clinical impression List by prescription"
#    (populated in C -> D -> B order to match the original authoring)
# ---------------------------------------------------------------------------
Set-QueryCell "C11" $s53 $true
Set-QueryCell "D11" $s54 $false
Set-QueryCell "B11" $s55 $true
$ws.Rows.Item(11).RowHeight = 60

# ---------------------------------------------------------------------------
# 6. Update the view's active selection to match the target.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C8").Select()
